$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("HARSHAD NAGTILAK", 22913.3785),
    @("HRUTWIK GARDI", 592.861),
    @("PRATIK RAUL", 29943.4216),
    @("PRATIK SHIRBHATE", 1752.8626499999998),
    @("SARANG THAKREY", 557.0264),
    @("SHUBHAM MUNDADA", 89968.73989999996),
    @("YUKTA SONIGRA", 16846.567150000003)
)

$startRow = 23
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = "2025-04"
}
